$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete MuSCs-sending block (old rows 8-10); the new TPM data
# only has 6 data rows (FAPs- and MuSCs-sending), not 9.
$ws.Rows.Item(10).EntireRow.Delete() | Out-Null
$ws.Rows.Item(9).EntireRow.Delete() | Out-Null
$ws.Rows.Item(8).EntireRow.Delete() | Out-Null

# Overwrite rows 2-7 with the recomputed TPM values
# Row 2: FAPs -> ECs
$ws.Cells.Item(2, 1).Value2 = "FAPs"
$ws.Cells.Item(2, 2).Value2 = "Dlk1"
$ws.Cells.Item(2, 3).Value2 = "Notch4"
$ws.Cells.Item(2, 4).Value2 = "ECs"
$ws.Cells.Item(2, 5).Value2 = 3
$ws.Cells.Item(2, 6).Value2 = 1
$ws.Cells.Item(2, 7).Value2 = 2.101448
$ws.Cells.Item(2, 8).Value2 = 6.304344
$ws.Cells.Item(2, 9).Value2 = 0.5480341737688159
$ws.Cells.Item(2, 10).Value2 = 0.5480341737688159
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 12).Value2 = 1
$ws.Cells.Item(2, 13).Value2 = 31.61061466666667
$ws.Cells.Item(2, 14).Value2 = 94.831844
$ws.Cells.Item(2, 15).Value2 = 0.8860472269592234
$ws.Cells.Item(2, 16).Value2 = 0.8860472269592234
$ws.Cells.Item(2, 17).Value2 = 66.42806297003733
$ws.Cells.Item(2, 18).Value2 = 597.8525667303361
$ws.Cells.Item(2, 19).Value2 = 0.4855841599467485
$ws.Cells.Item(2, 20).Value2 = 0.4855841599467485

# Row 3: FAPs -> FAPs
$ws.Cells.Item(3, 1).Value2 = "FAPs"
$ws.Cells.Item(3, 2).Value2 = "Dlk1"
$ws.Cells.Item(3, 3).Value2 = "Notch4"
$ws.Cells.Item(3, 4).Value2 = "FAPs"
$ws.Cells.Item(3, 5).Value2 = 3
$ws.Cells.Item(3, 6).Value2 = 1
$ws.Cells.Item(3, 7).Value2 = 2.101448
$ws.Cells.Item(3, 8).Value2 = 6.304344
$ws.Cells.Item(3, 9).Value2 = 0.5480341737688159
$ws.Cells.Item(3, 10).Value2 = 0.5480341737688159
$ws.Cells.Item(3, 11).Value2 = 3
$ws.Cells.Item(3, 12).Value2 = 1
$ws.Cells.Item(3, 13).Value2 = 1.672785333333334
$ws.Cells.Item(3, 14).Value2 = 5.018356000000001
$ws.Cells.Item(3, 15).Value2 = 0.04688826274109129
$ws.Cells.Item(3, 16).Value2 = 0.04688826274109129
$ws.Cells.Item(3, 17).Value2 = 3.515271393162667
$ws.Cells.Item(3, 18).Value2 = 31.63744253846401
$ws.Cells.Item(3, 19).Value2 = 0.02569637033076912
$ws.Cells.Item(3, 20).Value2 = 0.02569637033076912

# Row 4: FAPs -> MuSCs
$ws.Cells.Item(4, 1).Value2 = "FAPs"
$ws.Cells.Item(4, 2).Value2 = "Dlk1"
$ws.Cells.Item(4, 3).Value2 = "Notch4"
$ws.Cells.Item(4, 4).Value2 = "MuSCs"
$ws.Cells.Item(4, 5).Value2 = 3
$ws.Cells.Item(4, 6).Value2 = 1
$ws.Cells.Item(4, 7).Value2 = 2.101448
$ws.Cells.Item(4, 8).Value2 = 6.304344
$ws.Cells.Item(4, 9).Value2 = 0.5480341737688159
$ws.Cells.Item(4, 10).Value2 = 0.5480341737688159
$ws.Cells.Item(4, 11).Value2 = 3
$ws.Cells.Item(4, 12).Value2 = 1
$ws.Cells.Item(4, 13).Value2 = 2.392593
$ws.Cells.Item(4, 14).Value2 = 7.177778999999999
$ws.Cells.Item(4, 15).Value2 = 0.06706451029968528
$ws.Cells.Item(4, 16).Value2 = 0.06706451029968527
$ws.Cells.Item(4, 17).Value2 = 5.027909774664
$ws.Cells.Item(4, 18).Value2 = 45.25118797197599
$ws.Cells.Item(4, 19).Value2 = 0.03675364349129827
$ws.Cells.Item(4, 20).Value2 = 0.03675364349129826

# Row 5: MuSCs -> ECs
$ws.Cells.Item(5, 1).Value2 = "MuSCs"
$ws.Cells.Item(5, 2).Value2 = "Dlk1"
$ws.Cells.Item(5, 3).Value2 = "Notch4"
$ws.Cells.Item(5, 4).Value2 = "ECs"
$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 6).Value2 = 1
$ws.Cells.Item(5, 7).Value2 = 1.733072
$ws.Cells.Item(5, 8).Value2 = 5.199216
$ws.Cells.Item(5, 9).Value2 = 0.4519658262311841
$ws.Cells.Item(5, 10).Value2 = 0.4519658262311841
$ws.Cells.Item(5, 11).Value2 = 3
$ws.Cells.Item(5, 12).Value2 = 1
$ws.Cells.Item(5, 13).Value2 = 31.61061466666667
$ws.Cells.Item(5, 14).Value2 = 94.831844
$ws.Cells.Item(5, 15).Value2 = 0.8860472269592234
$ws.Cells.Item(5, 16).Value2 = 0.8860472269592234
$ws.Cells.Item(5, 17).Value2 = 54.78347118158933
$ws.Cells.Item(5, 18).Value2 = 493.051240634304
$ws.Cells.Item(5, 19).Value2 = 0.4004630670124749
$ws.Cells.Item(5, 20).Value2 = 0.4004630670124749

# Row 6: MuSCs -> FAPs
$ws.Cells.Item(6, 1).Value2 = "MuSCs"
$ws.Cells.Item(6, 2).Value2 = "Dlk1"
$ws.Cells.Item(6, 3).Value2 = "Notch4"
$ws.Cells.Item(6, 4).Value2 = "FAPs"
$ws.Cells.Item(6, 5).Value2 = 3
$ws.Cells.Item(6, 6).Value2 = 1
$ws.Cells.Item(6, 7).Value2 = 1.733072
$ws.Cells.Item(6, 8).Value2 = 5.199216
$ws.Cells.Item(6, 9).Value2 = 0.4519658262311841
$ws.Cells.Item(6, 10).Value2 = 0.4519658262311841
$ws.Cells.Item(6, 11).Value2 = 3
$ws.Cells.Item(6, 12).Value2 = 1
$ws.Cells.Item(6, 13).Value2 = 1.672785333333334
$ws.Cells.Item(6, 14).Value2 = 5.018356000000001
$ws.Cells.Item(6, 15).Value2 = 0.04688826274109129
$ws.Cells.Item(6, 16).Value2 = 0.04688826274109129
$ws.Cells.Item(6, 17).Value2 = 2.899057423210667
$ws.Cells.Item(6, 18).Value2 = 26.091516808896
$ws.Cells.Item(6, 19).Value2 = 0.02119189241032217
$ws.Cells.Item(6, 20).Value2 = 0.02119189241032217

# Row 7: MuSCs -> MuSCs
$ws.Cells.Item(7, 1).Value2 = "MuSCs"
$ws.Cells.Item(7, 2).Value2 = "Dlk1"
$ws.Cells.Item(7, 3).Value2 = "Notch4"
$ws.Cells.Item(7, 4).Value2 = "MuSCs"
$ws.Cells.Item(7, 5).Value2 = 3
$ws.Cells.Item(7, 6).Value2 = 1
$ws.Cells.Item(7, 7).Value2 = 1.733072
$ws.Cells.Item(7, 8).Value2 = 5.199216
$ws.Cells.Item(7, 9).Value2 = 0.4519658262311841
$ws.Cells.Item(7, 10).Value2 = 0.4519658262311841
$ws.Cells.Item(7, 11).Value2 = 3
$ws.Cells.Item(7, 12).Value2 = 1
$ws.Cells.Item(7, 13).Value2 = 2.392593
$ws.Cells.Item(7, 14).Value2 = 7.177778999999999
$ws.Cells.Item(7, 15).Value2 = 0.06706451029968528
$ws.Cells.Item(7, 16).Value2 = 0.06706451029968527
$ws.Cells.Item(7, 17).Value2 = 4.146535935696
$ws.Cells.Item(7, 18).Value2 = 37.31882342126399
$ws.Cells.Item(7, 19).Value2 = 0.03031086680838702
$ws.Cells.Item(7, 20).Value2 = 0.03031086680838701
